$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Ciboulette" sheet. It is
# inserted as row 431, pushing every existing row from 431 downward by one
# (431->432, ..., 462->463). Insert a fresh row at 431 first so the rows
# below shift down and keep their data/formatting intact.
$ws.Rows("431:431").Insert()

# Populate the newly inserted row with the new record's values.
$ws.Range("A431").Value = 9
$ws.Range("B431").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C431").Value = "Metropolitana"
$ws.Range("D431").Value = 44826
$ws.Range("E431").Value = 13
$ws.Range("F431").Value = 100112039
$ws.Range("G431").Value = "Ciboulette"
$ws.Range("H431").Value = "Sin especificar"
$ws.Range("I431").Value = "Primera"
$ws.Range("J431").Value = 430
$ws.Range("K431").Value = 1000
$ws.Range("L431").Value = 1000
$ws.Range("M431").Value = 1000
$ws.Range("N431").Value = "`$/docena de atados"
$ws.Range("O431").Value = "Región Metropolitana"
$ws.Range("P431").Value = 333
$ws.Range("Q431").Value = 3
$ws.Range("R431").Value = "Hortaliza"
